# Commit: "Added Necessary User Story / Location history requiring DB"
#
# 1. Append a brand-new 4-slide "feature" group (User Story, [User Accounts],
#    [Possible Technology], [Chosen Technology]) to the end of the deck by
#    duplicating the previous feature group (slides 55-58) and re-titling the
#    new user-story slide.
# 2. Re-stamp the "datetimeFigureOut" date placeholders (slide master + every
#    slide layout) from 5/11/2018 to 11/5/2018.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Duplicate the last feature group (User Story / User Accounts /
#    Possible Technology / Chosen Technology -- slides 55-58) and move the
#    four copies to the end of the deck, preserving their relative order.
# ---------------------------------------------------------------------------
$storySlideIndex = 55
$accountsSlideIndex = 56
$technologySlideIndex = 57
$chosenSlideIndex = 58

$dup1 = $p.Slides.Item($storySlideIndex).Duplicate()
$dup1.Item(1).MoveTo($p.Slides.Count)

$dup2 = $p.Slides.Item($accountsSlideIndex).Duplicate()
$dup2.Item(1).MoveTo($p.Slides.Count)

$dup3 = $p.Slides.Item($technologySlideIndex).Duplicate()
$dup3.Item(1).MoveTo($p.Slides.Count)

$dup4 = $p.Slides.Item($chosenSlideIndex).Duplicate()
$dup4.Item(1).MoveTo($p.Slides.Count)

# The new user-story slide is now the last-but-three slide; update its title
# text to the new user story.
$newStorySlide = $p.Slides.Item($p.Slides.Count - 3)
$newStorySlide.Shapes.Item(1).TextFrame.TextRange.Text = "As Kara the dedicated user, I want to save my landmarks online so I can access them from another device."

# ---------------------------------------------------------------------------
# 2) Fix up the "Last modified" date placeholder text everywhere it appears
#    (slide master + all slide layouts).
# ---------------------------------------------------------------------------
$oldDate = "5/11/2018"
$newDate = "11/5/2018"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "*Date*" -and $shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "*Date*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
